# REVER_DailyTracker: add Soniya's Oct 2-6 entries to Ram's sheet, update
# the WIP/Completed legend swatches on Ram's sheet, and move the active
# sheet/selection from Nirmal to Ram.

$wb = $excel.ActiveWorkbook

$wsRam    = $wb.Worksheets.Item("Ram")
$wsNirmal = $wb.Worksheets.Item("Nirmal")

# ---------------------------------------------------------------------
# 1) New data rows 4-8 on Ram's sheet.
#    Rows 4-6 (Holiday / Week off placeholders) share Nirmal's row-3
#    styling; rows 7-8 (Soniya's task) share Nirmal's row 6/7 styling,
#    with column F swapped for the new centred WIP/Completed look.
# ---------------------------------------------------------------------

$wsRam.Rows.Item(4).ClearFormats()
$wsNirmal.Range("A3:G3").Copy()
$wsRam.Range("A4:G4").PasteSpecial(-4122)

$wsRam.Rows.Item(5).ClearFormats()
$wsNirmal.Range("A4:G4").Copy()
$wsRam.Range("A5:G5").PasteSpecial(-4122)

$wsRam.Rows.Item(6).ClearFormats()
$wsNirmal.Range("A5:G5").Copy()
$wsRam.Range("A6:G6").PasteSpecial(-4122)

$wsRam.Rows.Item(7).ClearFormats()
$wsNirmal.Range("A6:G6").Copy()
$wsRam.Range("A7:G7").PasteSpecial(-4122)

$wsRam.Rows.Item(8).ClearFormats()
$wsNirmal.Range("A7:G7").Copy()
$wsRam.Range("A8:G8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Values for rows 4-6 (Holiday, Week off, Week off)
$wsRam.Range("A4").Value = 2
$wsRam.Range("B4").Value = 44106
$wsRam.Range("D4").Value = "Holiday"

$wsRam.Range("A5").Value = 3
$wsRam.Range("B5").Value = 44107
$wsRam.Range("D5").Value = "Week off"

$wsRam.Range("A6").Value = 4
$wsRam.Range("B6").Value = 44108
$wsRam.Range("D6").Value = "Week off"

# Values for rows 7-8 (Soniya / creating Setup file)
$wsRam.Range("A7").Value = 5
$wsRam.Range("B7").Value = 44109
$wsRam.Range("C7").Value = "Soniya"
$wsRam.Range("D7").Value = "creating Setup file"
$wsRam.Range("E7").Value = 0.5
$wsRam.Range("F7").Value = "WIP"
$wsRam.Range("F7").HorizontalAlignment = -4108

$wsRam.Range("A8").Value = 6
$wsRam.Range("B8").Value = 44110
$wsRam.Range("C8").Value = "Soniya"
$wsRam.Range("D8").Value = "creating Setup file"
$wsRam.Range("E8").Value = 1
$wsRam.Range("F8").Value = "Completed"
$wsRam.Range("F8").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 2) Legend swatches (rows 23 & 25) on Ram's sheet pick up WIP/Completed
#    labels, centred, matching the colours already used in the legend.
# ---------------------------------------------------------------------
$wsRam.Range("B23").Value = "WIP"
$wsRam.Range("B23").HorizontalAlignment = -4108

$wsRam.Range("B25").Value = "Completed"
$wsRam.Range("B25").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3) Selection / active-sheet bookkeeping: Nirmal gives up the active
#    tab (selecting its row 3) and Ram becomes the active sheet with
#    D14:D15 selected.
# ---------------------------------------------------------------------
$wsNirmal.Rows.Item(3).Select()

$wsRam.Activate()
$wsRam.Range("D14:D15").Select()
